$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.693.46'
$ws.Range("E2").Value = '  -3.43%  '
$ws.Range("D3").Value = '3.419.99'
$ws.Range("E3").Value = '  -5.23%  '
$ws.Range("E4").Value = '  +0.12%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '182.92'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  -10.68%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '533.19'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -5.95%  '
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.616'
$c.Style = 'Normal'
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("D8").Value = '3.412.44'
$ws.Range("E8").Value = '  -5.28%  '
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.631'
$c.Style = 'Normal'
$ws.Range("E10").Value = '  -6.37%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '58.09'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  -4.80%  '
$ws.Range("E12").Value = '  -10.45%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.0000257'
$c.Style = 'Normal'
$ws.Range("E13").Value = '  -10.68%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '9.44'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  -5.83%  '
$ws.Range("D15").Value = '3.965.30'
$ws.Range("E15").Value = '  -5.52%  '
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").Value = '3.414.07'
$ws.Range("E17").Value = '  -5.34%  '
$ws.Range("D18").Value = '65.528.14'
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '17.72'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  -6.03%  '
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '11.34'
$c.Style = 'Normal'
$ws.Range("E20").Value = '  -8.34%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '0.990'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -7.90%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '380.59'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  -5.39%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '83.80'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  -1.85%  '
$ws.Range("E24").Value = '  -9.20%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '10.92'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -16.64%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '11.67'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  -7.63%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '3.64'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  -8.13%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '2.69'
$c.Style = 'Normal'
$ws.Range("E28").Value = '  -8.71%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '8.60'
$c.Style = 'Normal'
$ws.Range("E29").Value = '  -8.82%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '687.72'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  +1.79%  '
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '30.11'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  -4.84%  '
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '6.79'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -19.52%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '11.29'
$c.Style = 'Normal'
$ws.Range("E33").Value = '  -7.62%  '
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '61.75'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  -3.36%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '0.106'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  -6.57%  '
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range("E36").Value = '  -0.05%  '
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '37.03'
$c.Style = 'Normal'
$ws.Range("E37").Value = '  -12.35%  '
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.390'
$c.Style = 'Normal'
$ws.Range("E38").Value = '  -7.62%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  -5.98%  '
$ws.Range("D41").Value = '2.901.58'
$ws.Range("E41").Value = '  -11.46%  '
$ws.Range("E42").Value = '  -12.64%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '2.69'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").Value = '0.0₃0632'
$ws.Range("E44").Value = '  -17.91%  '
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.0393'
$c.Style = 'Normal'
$ws.Range("E45").Value = '  -6.25%  '
$ws.Range("E46").Value = '  -15.26%  '
$ws.Range("E47").Value = '  -3.47%  '
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '134.80'
$c.Style = 'Normal'
$ws.Range("E48").Value = '  -3.63%  '
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '2.84'
$c.Style = 'Normal'
$ws.Range("E49").Value = '  -7.11%  '
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '2.58'
$c.Style = 'Normal'
$ws.Range("E50").Value = '  -5.64%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '7.73'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  -12.54%  '
